{"js": "// Apply the two edits described by the diff:\n//  1. Insert a new \"Meta description: ...\" paragraph right after the\n//     document's opening Heading1 title paragraph. \"Meta description\" is\n//     bold, the rest of the sentence is regular text.\n//  2. Near the end of the document, delete the duplicate bold\n//     \"Play Free Disco Diamonds Slot - A Retro Disco Experience\" paragraph\n//     and replace the text of the following italic paragraph (originally\n//     \"Read our review of the Disco Diamonds slot game...\") with the new\n//     image-generation prompt text, keeping the italic formatting intact.\n\nconst boldLabel = \"Meta description\";\nconst metaRest =\n  \": Read our review of the Disco Diamonds slot game. Dance to the rhythm with two special modes that offer payouts, free spins, and multipliers. Play for free now.\";\nconst promptText =\n  \"Prompt: Create a cartoon-style feature image that showcases the fun and unique identity of the \\\"Disco Diamonds\\\" game. The image should focus on a happy Maya warrior character wearing glasses, as a nod to the game's disco party theme. The image should be bright and colorful with funky disco elements, such as a disco ball and neon lights. The Maya warrior character should be drawn with a big smile and wearing stylish glasses, emphasizing the fun, upbeat feel of the game. The image should also include text that reads \\\"Disco Diamonds\\\", using a bold and playful font. The text should stand out and capture the attention of potential players. Overall, the feature image should capture the essence of the game's party theme while also showcasing its simple and engaging gameplay. The colorful and eye-catching design should entice players to give the game a try and experience the excitement of Disco Diamonds.\";\n\nconst body = context.document.body;\n\n// --- Edit 1: insert the Meta description paragraph after the title ---\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\n\nconst titlePara = paras.items[0];\nconst metaPara = titlePara.insertParagraph(boldLabel + metaRest, Word.InsertLocation.after);\n// The new paragraph otherwise inherits the Heading1 style from the title;\n// reset it back to the default body style.\nmetaPara.style = \"Normal\";\nawait context.sync();\n\n// Bold just the \"Meta description\" label, leaving the rest of the\n// sentence in regular formatting.\nconst labelRanges = metaPara.search(boldLabel, { matchCase: true });\nlabelRanges.load(\"items\");\nawait context.sync();\nlabelRanges.items[0].font.bold = true;\nawait context.sync();\n\n// --- Edit 2: drop the duplicate bold title paragraph near the bottom and\n//     update the following italic paragraph's text ---\nconst paras2 = body.paragraphs;\nparas2.load(\"items\");\nawait context.sync();\n\nconst count = paras2.items.length;\nconst duplicateTitlePara = paras2.items[count - 2];\nduplicateTitlePara.delete();\nawait context.sync();\n\nconst paras3 = body.paragraphs;\nparas3.load(\"items\");\nawait context.sync();\n\nconst count3 = paras3.items.length;\nconst descriptionPara = paras3.items[count3 - 1];\ndescriptionPara.insertText(promptText, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Apply the two edits described by the diff:\n#  1. Insert a new \"Meta description: ...\" paragraph right after the\n#     document's opening Heading1 title paragraph. \"Meta description\" is\n#     bold, the rest of the sentence is regular text.\n#  2. Near the end of the document, delete the duplicate bold\n#     \"Play Free Disco Diamonds Slot - A Retro Disco Experience\" paragraph\n#     and replace the text of the following italic paragraph (originally\n#     \"Read our review of the Disco Diamonds slot game...\") with the new\n#     image-generation prompt text, keeping the italic formatting intact.\n\n$d = $word.ActiveDocument\n\n$boldLabel = \"Meta description\"\n$metaRest = \": Read our review of the Disco Diamonds slot game. Dance to the rhythm with two special modes that offer payouts, free spins, and multipliers. Play for free now.\"\n$promptText = \"Prompt: Create a cartoon-style feature image that showcases the fun and unique identity of the `\"Disco Diamonds`\" game. The image should focus on a happy Maya warrior character wearing glasses, as a nod to the game's disco party theme. The image should be bright and colorful with funky disco elements, such as a disco ball and neon lights. The Maya warrior character should be drawn with a big smile and wearing stylish glasses, emphasizing the fun, upbeat feel of the game. The image should also include text that reads `\"Disco Diamonds`\", using a bold and playful font. The text should stand out and capture the attention of potential players. Overall, the feature image should capture the essence of the game's party theme while also showcasing its simple and engaging gameplay. The colorful and eye-catching design should entice players to give the game a try and experience the excitement of Disco Diamonds.\"\n\n# --- Edit 1: insert the Meta description paragraph after the title ---\n$titlePara = $d.Paragraphs.Item(1)\n$titlePara.Range.InsertParagraphAfter()\n\n$metaPara = $d.Paragraphs.Item(2)\n# The new paragraph otherwise inherits the Heading1 style from the title;\n# reset it back to the default body style.\n$metaPara.Style = \"Normal\"\n\n$metaRange = $metaPara.Range\n$insertStart = $metaRange.Start\n$metaRange.InsertAfter($boldLabel + $metaRest)\n\n# Bold just the \"Meta description\" label, leaving the rest of the\n# sentence in regular formatting.\n$boldRange = $d.Range($insertStart, $insertStart + $boldLabel.Length)\n$boldRange.Font.Bold = $true\n\n# --- Edit 2: drop the duplicate bold title paragraph near the bottom and\n#     update the following italic paragraph's text ---\n$count = $d.Paragraphs.Count\n$duplicateTitlePara = $d.Paragraphs.Item($count - 1)\n$duplicateTitlePara.Range.Delete()\n\n$newCount = $d.Paragraphs.Count\n$descriptionPara = $d.Paragraphs.Item($newCount)\n$fullRange = $descriptionPara.Range\n$textOnlyRange = $d.Range($fullRange.Start, $fullRange.End - 1)\n$textOnlyRange.Text = $promptText\n"}
